$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '55.970.46'
$ws.Range("E2").Value = '  +4.11%  '

# Row 3
$ws.Range("D3").Value = '2.505.38'
$ws.Range("E3").Value = '  +5.48%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.14'
$ws.Range("E5").Value = '  +7.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.37'
$ws.Range("E6").Value = '  +13.19%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").Value = '  +6.99%  '

# Row 9
$ws.Range("D9").Value = '2.522.57'
$ws.Range("E9").Value = '  +5.81%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0980'
$ws.Range("E10").Value = '  +4.81%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.63'
$ws.Range("E11").Value = '  +6.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.334'
$ws.Range("E12").Value = '  +6.05%  '

# Row 13
$ws.Range("E13").Value = '  +1.74%  '

# Row 14
$ws.Range("D14").Value = '2.936.51'
$ws.Range("E14").Value = '  +5.46%  '

# Row 15
$ws.Range("D15").Value = '55.964.86'
$ws.Range("E15").Value = '  +4.54%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.15'
$ws.Range("E16").Value = '  +9.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000136'
$ws.Range("E17").Value = '  +7.81%  '

# Row 18
$ws.Range("D18").Value = '2.519.61'
$ws.Range("E18").Value = '  +5.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.43'
$ws.Range("E19").Value = '  +6.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.24'
$ws.Range("E20").Value = '  +10.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.86'
$ws.Range("E21").Value = '  +4.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.82'
$ws.Range("E23").Value = '  +9.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.44'
$ws.Range("E24").Value = '  +4.75%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.411'
$ws.Range("E25").Value = '  +9.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  +11.29%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.72%  '

# Row 28
$ws.Range("D28").Value = '2.621.18'
$ws.Range("E28").Value = '  +5.91%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.47'
$ws.Range("E29").Value = '  +6.31%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0788'
$ws.Range("E30").Value = '  +13.30%  '

# Row 31
$ws.Range("E31").Value = '  +0.35%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '148.66'
$ws.Range("E32").Value = '  +0.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.41'
$ws.Range("E33").Value = '  +5.18%  '

# Row 34
$ws.Range("E34").Value = '  +8.54%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.24'
$ws.Range("E35").Value = '  +5.78%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.15'
$ws.Range("E36").Value = '  +11.43%  '

# Row 37
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.70'
$ws.Range("E37").Value = '  +8.91%  '

# Row 38
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.867'
$ws.Range("E38").Value = '  +11.51%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.22'
$ws.Range("E39").Value = '  +2.55%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.54'
$ws.Range("E40").Value = '  +8.74%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.618'
$ws.Range("E41").Value = '  +3.48%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0557'
$ws.Range("E42").Value = '  +6.71%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.993'
$ws.Range("E43").Value = '  -0.01%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.32'
$ws.Range("E44").Value = '  +8.97%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.83'
$ws.Range("E45").Value = '  +15.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '265.33'
$ws.Range("E46").Value = '  +27.22%  '

# Row 47
$ws.Range("E47").Value = '  +0.54%  '

# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0228'
$ws.Range("E48").Value = '  +6.07%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0909'
$ws.Range("E49").Value = '  +6.05%  '

# Row 50
$ws.Range("D50").Value = '1.956.11'
$ws.Range("E50").Value = '  +1.10%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.70'
$ws.Range("E51").Value = '  +8.28%  '
